$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "50-1="
$t.Cell(1, 2).Range.Text = "20+23="
$t.Cell(1, 3).Range.Text = "38+28="
$t.Cell(1, 4).Range.Text = "17+82="
$t.Cell(1, 5).Range.Text = "97-55="
$t.Cell(2, 1).Range.Text = "11+22="
$t.Cell(2, 2).Range.Text = "83-61="
$t.Cell(2, 3).Range.Text = "34-19="
$t.Cell(2, 4).Range.Text = "89-10="
$t.Cell(2, 5).Range.Text = "79+2="
$t.Cell(3, 1).Range.Text = "46+49="
$t.Cell(3, 2).Range.Text = "70-48="
$t.Cell(3, 3).Range.Text = "90-14="
$t.Cell(3, 4).Range.Text = "46+45="
$t.Cell(3, 5).Range.Text = "57-16="
$t.Cell(4, 1).Range.Text = "8+53="
$t.Cell(4, 2).Range.Text = "81-24="
$t.Cell(4, 3).Range.Text = "38-22="
$t.Cell(4, 4).Range.Text = "64-39="
$t.Cell(4, 5).Range.Text = "32-22="
$t.Cell(5, 1).Range.Text = "10+11="
$t.Cell(5, 2).Range.Text = "53-16="
$t.Cell(5, 3).Range.Text = "49-36="
$t.Cell(5, 4).Range.Text = "17+17="
$t.Cell(5, 5).Range.Text = "98-57="
$t.Cell(6, 1).Range.Text = "23+50="
$t.Cell(6, 2).Range.Text = "62-0="
$t.Cell(6, 3).Range.Text = "70-45="
$t.Cell(6, 4).Range.Text = "83-36="
$t.Cell(6, 5).Range.Text = "17-4="
$t.Cell(7, 1).Range.Text = "72+8="
$t.Cell(7, 2).Range.Text = "62-62="
$t.Cell(7, 3).Range.Text = "92-16="
$t.Cell(7, 4).Range.Text = "50+42="
$t.Cell(7, 5).Range.Text = "16+65="
$t.Cell(8, 1).Range.Text = "99-98="
$t.Cell(8, 2).Range.Text = "25+29="
$t.Cell(8, 3).Range.Text = "97-34="
$t.Cell(8, 4).Range.Text = "23+2="
$t.Cell(8, 5).Range.Text = "91-79="
$t.Cell(9, 1).Range.Text = "73+17="
$t.Cell(9, 2).Range.Text = "75-58="
$t.Cell(9, 3).Range.Text = "65-42="
$t.Cell(9, 4).Range.Text = "19+32="
$t.Cell(9, 5).Range.Text = "14+72="
$t.Cell(10, 1).Range.Text = "81-74="
$t.Cell(10, 2).Range.Text = "75+19="
$t.Cell(10, 3).Range.Text = "10+66="
$t.Cell(10, 4).Range.Text = "9+19="
$t.Cell(10, 5).Range.Text = "21-10="
$t.Cell(11, 1).Range.Text = "20+29="
$t.Cell(11, 2).Range.Text = "63-48="
$t.Cell(11, 3).Range.Text = "78-53="
$t.Cell(11, 4).Range.Text = "82-3="
$t.Cell(11, 5).Range.Text = "5+57="
$t.Cell(12, 1).Range.Text = "76+10="
$t.Cell(12, 2).Range.Text = "54-8="
$t.Cell(12, 3).Range.Text = "77-2="
$t.Cell(12, 4).Range.Text = "82-19="
$t.Cell(12, 5).Range.Text = "40+37="
$t.Cell(13, 1).Range.Text = "14+52="
$t.Cell(13, 2).Range.Text = "63-31="
$t.Cell(13, 3).Range.Text = "93-23="
$t.Cell(13, 4).Range.Text = "78-76="
$t.Cell(13, 5).Range.Text = "87-79="
$t.Cell(14, 1).Range.Text = "20+66="
$t.Cell(14, 2).Range.Text = "76-46="
$t.Cell(14, 3).Range.Text = "76-61="
$t.Cell(14, 4).Range.Text = "90-16="
$t.Cell(14, 5).Range.Text = "62-4="
$t.Cell(15, 1).Range.Text = "85-41="
$t.Cell(15, 2).Range.Text = "40-7="
$t.Cell(15, 3).Range.Text = "73-64="
$t.Cell(15, 4).Range.Text = "7+6="
$t.Cell(15, 5).Range.Text = "32+19="
$t.Cell(16, 1).Range.Text = "74+8="
$t.Cell(16, 2).Range.Text = "21+33="
$t.Cell(16, 3).Range.Text = "0+78="
$t.Cell(16, 4).Range.Text = "58+16="
$t.Cell(16, 5).Range.Text = "55+3="
$t.Cell(17, 1).Range.Text = "9+81="
$t.Cell(17, 2).Range.Text = "51+21="
$t.Cell(17, 3).Range.Text = "51+44="
$t.Cell(17, 4).Range.Text = "43+27="
$t.Cell(17, 5).Range.Text = "11-10="
$t.Cell(18, 1).Range.Text = "27-19="
$t.Cell(18, 2).Range.Text = "78-9="
$t.Cell(18, 3).Range.Text = "26+53="
$t.Cell(18, 4).Range.Text = "15+66="
$t.Cell(18, 5).Range.Text = "55+4="
$t.Cell(19, 1).Range.Text = "97-45="
$t.Cell(19, 2).Range.Text = "31-22="
$t.Cell(19, 3).Range.Text = "99-4="
$t.Cell(19, 4).Range.Text = "84-29="
$t.Cell(19, 5).Range.Text = "60+8="
$t.Cell(20, 1).Range.Text = "38+26="
$t.Cell(20, 2).Range.Text = "32+43="
$t.Cell(20, 3).Range.Text = "75-21="
$t.Cell(20, 4).Range.Text = "94-28="
$t.Cell(20, 5).Range.Text = "35+1="
